$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The averaged-intensities table lists one "scheme" per row (col A = numeric
# scheme id, col B = scheme name). Re-running the notebook with the new
# spiral sampling schemes re-generated the scheme list: Gaussian-Quadrature
# and the three new Spiral-* schemes were inserted into the master scheme
# order right after the single/ring schemes, pushing the rotation/hex-grid
# schemes further down, and three brand-new rows were appended at the
# bottom of the sheet for the schemes that no longer fit in the original
# 0-14 row range.

# Rows 10-16: relabel to match the regenerated scheme order.
$ws.Cells.Item(10, 2).Value = "Gaussian-Quadrature"
$ws.Cells.Item(11, 2).Value = "Spiral-90deg-10rot-5space"
$ws.Cells.Item(12, 2).Value = "Spiral-90deg-15rot-5space"
$ws.Cells.Item(13, 2).Value = "Spiral-90deg-10rot-3space"
$ws.Cells.Item(14, 2).Value = "NoRotation-tilt60deg"
$ws.Cells.Item(15, 2).Value = "Rotation-NoTilt"
$ws.Cells.Item(16, 2).Value = "Rotation-60detTilt"

# New rows 17-19, continuing the scheme id sequence (15, 16, 17) and the
# regenerated scheme-name order (the hex-grid schemes that were displaced
# from rows 13-15).
$ws.Cells.Item(17, 1).Value = 15
$ws.Cells.Item(17, 2).Value = "HexGrid-90degTilt5degRes"
for ($c = 3; $c -le 16; $c++) {
    $ws.Cells.Item(17, $c).Value = 1
}

$ws.Cells.Item(18, 1).Value = 16
$ws.Cells.Item(18, 2).Value = "HexGrid-90degTilt22p5degRes"
for ($c = 3; $c -le 16; $c++) {
    $ws.Cells.Item(18, $c).Value = 1
}

$ws.Cells.Item(19, 1).Value = 17
$ws.Cells.Item(19, 2).Value = "HexGrid-60degTilt5degRes"
for ($c = 3; $c -le 16; $c++) {
    $ws.Cells.Item(19, $c).Value = 1
}

# Carry the bordered/bold/centered "index" style from column A of the last
# existing row down onto the newly added rows (matches A2:A16 formatting).
$ws.Range("A16").Copy()
$ws.Range("A17:A19").PasteSpecial(-4122)
$excel.CutCopyMode = $false
